$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road | Potion
$ws.Range("H17").Value = 2666.068
$ws.Range("J17").Value = 2725.5117
$ws.Range("L17").Value = 8176.5351
$ws.Range("N17").Value = -8512.535100000001

# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 1551.25
$ws.Range("I40").Value = 966.6667
$ws.Range("K40").Value = 966.6667
$ws.Range("M40").Value = -791.6667

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4534.481
$ws.Range("I138").Value = 1751.5454
$ws.Range("J138").Value = 6575.3
$ws.Range("K138").Value = 5254.6362
$ws.Range("L138").Value = 19725.9
$ws.Range("M138").Value = -114.6361999999999
$ws.Range("N138").Value = -30005.9

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 5120.6816
$ws.Range("I141").Value = 3457
$ws.Range("K141").Value = 10371
$ws.Range("M141").Value = -5191

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 16664.07
$ws.Range("I32").Value = 18317.049
$ws.Range("J32").Value = 5276.8887
$ws.Range("K32").Value = 18317.049
$ws.Range("L32").Value = 5276.8887
$ws.Range("M32").Value = -18030.049
$ws.Range("N32").Value = -5850.8887

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 2047.0883
$ws.Range("I74").Value = 2034.2333
$ws.Range("K74").Value = 2034.2333
$ws.Range("M74").Value = -1160.2333

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 2047.0883
$ws.Range("I77").Value = 2034.2333
$ws.Range("K77").Value = 10171.1665
$ws.Range("M77").Value = -5803.166500000001

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 50000
$ws.Range("I88").Value = 50000
$ws.Range("K88").Value = 50000
$ws.Range("M88").Value = -49594

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 50000
$ws.Range("I91").Value = 50000
$ws.Range("K91").Value = 50000
$ws.Range("M91").Value = -48596

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 2692.842
$ws.Range("I110").Value = 1311.4
$ws.Range("J110").Value = 4227.778
$ws.Range("K110").Value = 1311.4
$ws.Range("L110").Value = 4227.778
$ws.Range("M110").Value = 733.5999999999999
$ws.Range("N110").Value = -8317.778

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2328.8386
$ws.Range("I132").Value = 1732.1177
$ws.Range("J132").Value = 3053.4285
$ws.Range("K132").Value = 5196.3531
$ws.Range("L132").Value = 9160.2855
$ws.Range("M132").Value = -2666.3531
$ws.Range("N132").Value = -14220.2855

$ws = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight | Mythrite Nugget
$ws.Range("H64").Value = 554.1111
$ws.Range("J64").Value = 581.1667
$ws.Range("L64").Value = 581.1667
$ws.Range("N64").Value = -1031.1667

# Row 67: Bearing the Brunt (L) | Mythrite Nugget
$ws.Range("H67").Value = 554.1111
$ws.Range("J67").Value = 581.1667
$ws.Range("L67").Value = 581.1667
$ws.Range("N67").Value = -2141.1667

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 1808.772
$ws.Range("I86").Value = 1796.0944
$ws.Range("K86").Value = 1796.0944
$ws.Range("M86").Value = -673.0944

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 1808.772
$ws.Range("I89").Value = 1796.0944
$ws.Range("K89").Value = 8980.472
$ws.Range("M89").Value = -3364.472

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 914.03845
$ws.Range("I107").Value = 879.45
$ws.Range("J107").Value = 1029.3334
$ws.Range("K107").Value = 879.45
$ws.Range("L107").Value = 1029.3334
$ws.Range("M107").Value = 1040.55
$ws.Range("N107").Value = -4869.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2652.8064
$ws.Range("I31").Value = 1871.1538
$ws.Range("J31").Value = 3217.3333
$ws.Range("K31").Value = 1871.1538
$ws.Range("L31").Value = 3217.3333
$ws.Range("M31").Value = -1576.1538
$ws.Range("N31").Value = -3807.3333

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2652.8064
$ws.Range("I34").Value = 1871.1538
$ws.Range("J34").Value = 3217.3333
$ws.Range("K34").Value = 1871.1538
$ws.Range("L34").Value = 3217.3333
$ws.Range("M34").Value = -1669.1538
$ws.Range("N34").Value = -3621.3333

# Row 62: Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2451.5908
$ws.Range("I132").Value = 2210.8484
$ws.Range("K132").Value = 6632.5452
$ws.Range("M132").Value = -4102.5452

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water | Boiled Egg
$ws.Range("H4").Value = 379.2
$ws.Range("I4").Value = 225
$ws.Range("J4").Value = 996
$ws.Range("K4").Value = 675
$ws.Range("L4").Value = 2988
$ws.Range("M4").Value = -563
$ws.Range("N4").Value = -3212

# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 11912342
$ws.Range("I5").Value = 420.2
$ws.Range("J5").Value = 41692144
$ws.Range("K5").Value = 1260.6
$ws.Range("L5").Value = 125076432
$ws.Range("M5").Value = -1148.6
$ws.Range("N5").Value = -125076656

# Row 14: Keep Your Powder Dry | Kukuru Powder
$ws.Range("H14").Value = 1466.381
$ws.Range("I14").Value = 1466.381
$ws.Range("K14").Value = 4399.143
$ws.Range("M14").Value = -4226.143

# Row 109: Cure for What Ails | Purple Carrot Juice
$ws.Range("H109").Value = 2160.5
$ws.Range("I109").Value = 1363.5
$ws.Range("J109").Value = 2957.5
$ws.Range("K109").Value = 4090.5
$ws.Range("L109").Value = 8872.5
$ws.Range("M109").Value = -3050.5
$ws.Range("N109").Value = -10952.5

# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 11912342
$ws.Range("I135").Value = 420.2
$ws.Range("J135").Value = 41692144
$ws.Range("K135").Value = 3781.8
$ws.Range("L135").Value = 375229296
$ws.Range("M135").Value = -1246.8
$ws.Range("N135").Value = -375234366

$ws = $wb.Worksheets.Item("GSM")
# Row 119: Bulking Up | Dwarven Mythril Rapier
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 11374.714
$ws.Range("I132").Value = 4741.75
$ws.Range("K132").Value = 14225.25
$ws.Range("M132").Value = -11695.25

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 6805.804
$ws.Range("I122").Value = 6257.3105
$ws.Range("J122").Value = 7741.4707
$ws.Range("K122").Value = 18771.9315
$ws.Range("L122").Value = 23224.4121
$ws.Range("M122").Value = -16321.9315
$ws.Range("N122").Value = -28124.4121

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3165.7632
$ws.Range("I132").Value = 2708.5334
$ws.Range("K132").Value = 8125.600199999999
$ws.Range("M132").Value = -5595.600199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 1596.8182
$ws.Range("I122").Value = 1070.625
$ws.Range("K122").Value = 3211.875
$ws.Range("M122").Value = -761.875
